$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 33
$ws.Range("I2").Value = 95
$ws.Range("J2").Value = 380
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 87
$ws.Range("M2").Value = 5
$ws.Range("N2").Value = 71
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 1
$ws.Range("R2").Value = 5
$ws.Range("S2").Value = 51
$ws.Range("T2").Value = 64
$ws.Range("U2").Value = 4
$ws.Range("V2").Value = 557
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 591
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 5
$ws.Range("AA2").Value = 3
